$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Missing sex" column (column E) entirely, shifting F (Age range) left to E,
# and what was G (unused after shift) disappears.
$ws.Range("E1:E4").Delete()
